# posts.xlsx update: remove the post row that was taken down
# ("「もし敬意が落ちたら、それと一緒に全部落ちるのだ」", row 593), causing all
# subsequent rows to shift up by one and the sheet's used range to shrink by
# one row (A1:C668 -> A1:C667).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the entire row 593 and shift the rows below it up.
$ws.Rows.Item(593).Delete()
